$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (and implicitly the <sheet> entry in workbook.xml)
$ws.Name = 'XI BC 1'

# Update the exam info header block (rows 1-2)
$ws.Range('B1').Value = 'SIMULASI PTS KELAS XI'
$ws.Range('E1').Value = 'Kamis 04 Maret 2021 09:00 - 10:00'
$ws.Range('B2').Value = '60 menit'

# Replace the list of student names (column A, rows 11-44) and
# append the two new students (rows 45-46)
$names = @(
    'AMIRAH NUR AINI',
    'ANGEL LYKA SARI',
    'Ayu Yuliana',
    'AINI PUTRI AZIZAH',
    'Alief Randhinka Putra',
    'Arya Juliyawan',
    'Aditya Rafly Fatahhudin',
    'ABDUL RAHMAN',
    'ANDI RADITA DAFA',
    'ANISA PURI APRILIA',
    'Chika Yuliani',
    'INTAN RAMADHANIA',
    'INTAN NURAINI',
    'DWI RIYANTO',
    'DINA NOVIA PUTRI',
    'Fikri Hawari',
    'FASHA RIANI PUTRI',
    'JUANDA',
    'DEFIRA ARYANTI',
    'MUTIARA PUTRI AZZAHRA',
    'NUR FARIZ DIRJA',
    'Hilda Al Kayis',
    'Iip Dany Budi Utomo',
    'NABILAH HARSIKA',
    'PITRIA',
    'Neng Silvi Aprilianti',
    'MUHAMMAD KHALIL FERGAL',
    'Muhammad Fadhil Nur Iskandar',
    'MUHAMAD KEVIN',
    'PUTRI SALSABILLA',
    'Rara Safira',
    'WINIE LAILATUL RHAMADANTI WIBOWO',
    'Sofi Dwiyanti',
    'Siska Audina Fadilah',
    'SYIFA FADILLAH',
    'Viona'
)

# New Nilai PG / Nilai Total scores for rows 11-46 (Nilai Esai stays 0)
$scores = @(
    85,
    66,
    71,
    68,
    84,
    71,
    75,
    76,
    74,
    66,
    75,
    91,
    55,
    55,
    57,
    67,
    91,
    71,
    82,
    75,
    87,
    77,
    91,
    84,
    85,
    88,
    80,
    68,
    69,
    54,
    91,
    57,
    85,
    58,
    73,
    64
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = 11 + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $scores[$i]
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = $scores[$i]
}

Write-Host ('Updated sheet ' + $ws.Name)
